$wb = $excel.ActiveWorkbook

# ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 432.66666  # H12: 499.4 -> 432.66666
$ws.Cells.Item(12, 9).Value = 432.66666  # I12: 499.4 -> 432.66666
$ws.Cells.Item(12, 11).Value = 432.66666  # K12: 499.4 -> 432.66666
$ws.Cells.Item(12, 13).Value = -262.66666  # M12: -329.4 -> -262.66666

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 428.25  # H33: 445.54544 -> 428.25
$ws.Cells.Item(33, 9).Value = 428.25  # I33: 445.54544 -> 428.25
$ws.Cells.Item(33, 11).Value = 428.25  # K33: 445.54544 -> 428.25
$ws.Cells.Item(33, 13).Value = -199.25  # M33: -216.54544 -> -199.25

# ALC row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(87, 8).Value = 89997.75  # H87: 89998.25 -> 89997.75
$ws.Cells.Item(87, 10).Value = 89997.75  # J87: 89998.25 -> 89997.75
$ws.Cells.Item(87, 12).Value = 89997.75  # L87: 89998.25 -> 89997.75
$ws.Cells.Item(87, 14).Value = -92493.75  # N87: -92494.25 -> -92493.75

# ALC row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(90, 8).Value = 89997.75  # H90: 89998.25 -> 89997.75
$ws.Cells.Item(90, 10).Value = 89997.75  # J90: 89998.25 -> 89997.75
$ws.Cells.Item(90, 12).Value = 269993.25  # L90: 269994.75 -> 269993.25
$ws.Cells.Item(90, 14).Value = -282473.25  # N90: -282474.75 -> -282473.25

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 459.125  # H135: 459.25 -> 459.125
$ws.Cells.Item(135, 9).Value = 453.35715  # I135: 453.5 -> 453.35715
$ws.Cells.Item(135, 11).Value = 4080.21435  # K135: 4081.5 -> 4080.21435
$ws.Cells.Item(135, 13).Value = -1545.21435  # M135: -1546.5 -> -1545.21435

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 3165.2  # H137: 2950.2273 -> 3165.2
$ws.Cells.Item(137, 9).Value = 2960  # I137: 2976.9167 -> 2960
$ws.Cells.Item(137, 10).Value = 3473  # J137: 2918.2 -> 3473
$ws.Cells.Item(137, 11).Value = 8880  # K137: 8930.750100000001 -> 8880
$ws.Cells.Item(137, 12).Value = 10419  # L137: 8754.599999999999 -> 10419
$ws.Cells.Item(137, 13).Value = -6330  # M137: -6380.750100000001 -> -6330
$ws.Cells.Item(137, 14).Value = -15519  # N137: -13854.6 -> -15519

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2337.1853  # H138: 2364.48 -> 2337.1853
$ws.Cells.Item(138, 9).Value = 1626.1538  # I138: 1678.5834 -> 1626.1538
$ws.Cells.Item(138, 10).Value = 2997.4285  # J138: 2997.6155 -> 2997.4285
$ws.Cells.Item(138, 11).Value = 4878.4614  # K138: 5035.7502 -> 4878.4614
$ws.Cells.Item(138, 12).Value = 8992.2855  # L138: 8992.8465 -> 8992.2855
$ws.Cells.Item(138, 13).Value = 261.5385999999999  # M138: 104.2497999999996 -> 261.5385999999999
$ws.Cells.Item(138, 14).Value = -19272.2855  # N138: -19272.8465 -> -19272.2855

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 276.5  # H5: 399.5 -> 276.5
$ws.Cells.Item(5, 9).Value = 180.71428  # I5: 249 -> 180.71428
$ws.Cells.Item(5, 10).Value = 500  # J5: 550 -> 500
$ws.Cells.Item(5, 11).Value = 180.71428  # K5: 249 -> 180.71428
$ws.Cells.Item(5, 12).Value = 500  # L5: 550 -> 500
$ws.Cells.Item(5, 13).Value = -68.71428  # M5: -137 -> -68.71428
$ws.Cells.Item(5, 14).Value = -724  # N5: -774 -> -724

# ARM row 46
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(46, 8).Value = 22025.6  # H46: 29046 -> 22025.6
$ws.Cells.Item(46, 10).Value = 14497.5  # J46: 17500 -> 14497.5
$ws.Cells.Item(46, 12).Value = 14497.5  # L46: 17500 -> 14497.5
$ws.Cells.Item(46, 14).Value = -15135.5  # N46: -18138 -> -15135.5

# ARM row 80
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(80, 8).Value = 26565.857  # H80: 26566.143 -> 26565.857
$ws.Cells.Item(80, 10).Value = 29992.4  # J80: 29992.8 -> 29992.4
$ws.Cells.Item(80, 12).Value = 29992.4  # L80: 29992.8 -> 29992.4
$ws.Cells.Item(80, 14).Value = -31988.4  # N80: -31988.8 -> -31988.4

# ARM row 83
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(83, 8).Value = 26565.857  # H83: 26566.143 -> 26565.857
$ws.Cells.Item(83, 10).Value = 29992.4  # J83: 29992.8 -> 29992.4
$ws.Cells.Item(83, 12).Value = 89977.20000000001  # L83: 89978.39999999999 -> 89977.20000000001
$ws.Cells.Item(83, 14).Value = -99961.20000000001  # N83: -99962.39999999999 -> -99961.20000000001

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 1775.7778  # H122: 2264.8333 -> 1775.7778
$ws.Cells.Item(122, 9).Value = 1873  # I122: 2264.8333 -> 1873
$ws.Cells.Item(122, 10).Value = 998  # J122: 0 -> 998
$ws.Cells.Item(122, 11).Value = 5619  # K122: 6794.499899999999 -> 5619
$ws.Cells.Item(122, 12).Value = 2994  # L122: 0 -> 2994
$ws.Cells.Item(122, 13).Value = -3169  # M122: -4344.499899999999 -> -3169
$ws.Cells.Item(122, 14).Value = -7894  # N122: None -> -7894

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 276.5  # H4: 399.5 -> 276.5
$ws.Cells.Item(4, 9).Value = 180.71428  # I4: 249 -> 180.71428
$ws.Cells.Item(4, 10).Value = 500  # J4: 550 -> 500
$ws.Cells.Item(4, 11).Value = 180.71428  # K4: 249 -> 180.71428
$ws.Cells.Item(4, 12).Value = 500  # L4: 550 -> 500
$ws.Cells.Item(4, 13).Value = -65.71428  # M4: -134 -> -65.71428
$ws.Cells.Item(4, 14).Value = -730  # N4: -780 -> -730

# BSM row 55
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(55, 8).Value = 0  # H55: 33333 -> 0
$ws.Cells.Item(55, 10).Value = 0  # J55: 33333 -> 0
$ws.Cells.Item(55, 12).ClearContents()  # L55: 33333 -> (removed)
$ws.Cells.Item(55, 14).Value = 0  # N55: -33879 -> 0

# BSM row 82
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 32863.547  # H82: 35150.3 -> 32863.547
$ws.Cells.Item(82, 9).Value = 13617.833  # I82: 14341.4 -> 13617.833
$ws.Cells.Item(82, 10).Value = 55958.4  # J82: 55959.2 -> 55958.4
$ws.Cells.Item(82, 11).Value = 13617.833  # K82: 14341.4 -> 13617.833
$ws.Cells.Item(82, 12).Value = 55958.4  # L82: 55959.2 -> 55958.4
$ws.Cells.Item(82, 13).Value = -13234.833  # M82: -13958.4 -> -13234.833
$ws.Cells.Item(82, 14).Value = -56724.4  # N82: -56725.2 -> -56724.4

# BSM row 85
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(85, 8).Value = 32863.547  # H85: 35150.3 -> 32863.547
$ws.Cells.Item(85, 9).Value = 13617.833  # I85: 14341.4 -> 13617.833
$ws.Cells.Item(85, 10).Value = 55958.4  # J85: 55959.2 -> 55958.4
$ws.Cells.Item(85, 11).Value = 13617.833  # K85: 14341.4 -> 13617.833
$ws.Cells.Item(85, 12).Value = 55958.4  # L85: 55959.2 -> 55958.4
$ws.Cells.Item(85, 13).Value = -12291.833  # M85: -13015.4 -> -12291.833
$ws.Cells.Item(85, 14).Value = -58610.4  # N85: -58611.2 -> -58610.4

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 68.21738999999999  # H7: 66.61905 -> 68.21738999999999
$ws.Cells.Item(7, 9).Value = 73.04761999999999  # I7: 66.61905 -> 73.04761999999999
$ws.Cells.Item(7, 10).Value = 17.5  # J7: 0 -> 17.5
$ws.Cells.Item(7, 11).Value = 73.04761999999999  # K7: 66.61905 -> 73.04761999999999
$ws.Cells.Item(7, 12).Value = 17.5  # L7: 0 -> 17.5
$ws.Cells.Item(7, 13).Value = 39.95238000000001  # M7: 46.38095 -> 39.95238000000001
$ws.Cells.Item(7, 14).Value = -243.5  # N7: None -> -243.5

# CRP row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(59, 8).Value = 41214.57  # H59: 42717 -> 41214.57
$ws.Cells.Item(59, 9).Value = 47752  # I59: 48752 -> 47752
$ws.Cells.Item(59, 10).Value = 38599.6  # J59: 39699.5 -> 38599.6
$ws.Cells.Item(59, 11).Value = 47752  # K59: 48752 -> 47752
$ws.Cells.Item(59, 12).Value = 38599.6  # L59: 39699.5 -> 38599.6
$ws.Cells.Item(59, 13).Value = -46607  # M59: -47607 -> -46607
$ws.Cells.Item(59, 14).Value = -40889.6  # N59: -41989.5 -> -40889.6

# CRP row 68
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(68, 8).Value = 49993.2  # H68: 49995.5 -> 49993.2
$ws.Cells.Item(68, 10).Value = 49993.2  # J68: 49995.5 -> 49993.2
$ws.Cells.Item(68, 12).Value = 49993.2  # L68: 49995.5 -> 49993.2
$ws.Cells.Item(68, 14).Value = -51491.2  # N68: -51493.5 -> -51491.2

# CRP row 71
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(71, 8).Value = 49993.2  # H71: 49995.5 -> 49993.2
$ws.Cells.Item(71, 10).Value = 49993.2  # J71: 49995.5 -> 49993.2
$ws.Cells.Item(71, 12).Value = 149979.6  # L71: 149986.5 -> 149979.6
$ws.Cells.Item(71, 14).Value = -157467.6  # N71: -157474.5 -> -157467.6

# CRP row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(74, 8).Value = 36342  # H74: 36286.2 -> 36342
$ws.Cells.Item(74, 10).Value = 36342  # J74: 36286.2 -> 36342
$ws.Cells.Item(74, 12).Value = 36342  # L74: 36286.2 -> 36342
$ws.Cells.Item(74, 14).Value = -38090  # N74: -38034.2 -> -38090

# CRP row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(77, 8).Value = 36342  # H77: 36286.2 -> 36342
$ws.Cells.Item(77, 10).Value = 36342  # J77: 36286.2 -> 36342
$ws.Cells.Item(77, 12).Value = 109026  # L77: 108858.6 -> 109026
$ws.Cells.Item(77, 14).Value = -117762  # N77: -117594.6 -> -117762

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 4456.1875  # H134: 4634.8667 -> 4456.1875
$ws.Cells.Item(134, 9).Value = 4286.6  # I134: 4465.9287 -> 4286.6
$ws.Cells.Item(134, 11).Value = 12859.8  # K134: 13397.7861 -> 12859.8
$ws.Cells.Item(134, 13).Value = -10324.8  # M134: -10862.7861 -> -10324.8

# CUL row 29
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(29, 8).Value = 263  # H29: 203.5 -> 263
$ws.Cells.Item(29, 10).Value = 19.5  # J29: 21.333334 -> 19.5
$ws.Cells.Item(29, 12).Value = 58.5  # L29: 64.00000199999999 -> 58.5
$ws.Cells.Item(29, 14).Value = -612.5  # N29: -618.000002 -> -612.5

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 120  # H2: 114.111115 -> 120
$ws.Cells.Item(2, 9).Value = 126.833336  # I2: 122.07692 -> 126.833336
$ws.Cells.Item(2, 10).Value = 99.5  # J2: 93.40000000000001 -> 99.5
$ws.Cells.Item(2, 11).Value = 126.833336  # K2: 122.07692 -> 126.833336
$ws.Cells.Item(2, 12).Value = 99.5  # L2: 93.40000000000001 -> 99.5
$ws.Cells.Item(2, 13).Value = -13.833336  # M2: -9.076920000000001 -> -13.833336
$ws.Cells.Item(2, 14).Value = -325.5  # N2: -319.4 -> -325.5

# GSM row 46
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 15562.728  # H46: 12119 -> 15562.728
$ws.Cells.Item(46, 10).Value = 50000  # J46: 0 -> 50000
$ws.Cells.Item(46, 12).Value = 50000  # L46: 0 -> 50000
$ws.Cells.Item(46, 14).Value = -50312  # N46: None -> -50312

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1791.1111  # H102: 1860 -> 1791.1111
$ws.Cells.Item(102, 9).Value = 1808.2354  # I102: 1881.1765 -> 1808.2354
$ws.Cells.Item(102, 11).Value = 1808.2354  # K102: 1881.1765 -> 1808.2354
$ws.Cells.Item(102, 13).Value = -186.2354  # M102: -259.1765 -> -186.2354

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 550  # H22: 525 -> 550
$ws.Cells.Item(22, 9).Value = 550  # I22: 525 -> 550
$ws.Cells.Item(22, 11).Value = 550  # K22: 525 -> 550
$ws.Cells.Item(22, 13).Value = -255  # M22: -230 -> -255

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 550  # H27: 525 -> 550
$ws.Cells.Item(27, 9).Value = 550  # I27: 525 -> 550
$ws.Cells.Item(27, 11).Value = 550  # K27: 525 -> 550
$ws.Cells.Item(27, 13).Value = -443  # M27: -418 -> -443

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1957.2333  # H46: 1980.5333 -> 1957.2333
$ws.Cells.Item(46, 9).Value = 1446  # I46: 1515.9 -> 1446
$ws.Cells.Item(46, 11).Value = 1446  # K46: 1515.9 -> 1446
$ws.Cells.Item(46, 13).Value = -1258  # M46: -1327.9 -> -1258

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1117.3125  # H55: 1138.4667 -> 1117.3125
$ws.Cells.Item(55, 10).Value = 1532.375  # J55: 1637 -> 1532.375
$ws.Cells.Item(55, 12).Value = 1532.375  # L55: 1637 -> 1532.375
$ws.Cells.Item(55, 14).Value = -1878.375  # N55: -1983 -> -1878.375

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 17972.5  # H122: 13381.333 -> 17972.5
$ws.Cells.Item(122, 9).Value = 7630  # I122: 6472.5 -> 7630
$ws.Cells.Item(122, 10).Value = 49000  # J122: 27199 -> 49000
$ws.Cells.Item(122, 11).Value = 22890  # K122: 19417.5 -> 22890
$ws.Cells.Item(122, 12).Value = 147000  # L122: 81597 -> 147000
$ws.Cells.Item(122, 13).Value = -20440  # M122: -16967.5 -> -20440
$ws.Cells.Item(122, 14).Value = -151900  # N122: -86497 -> -151900

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 1879.0588  # H132: 1833.1666 -> 1879.0588
$ws.Cells.Item(132, 9).Value = 1032  # I132: 1035 -> 1032
$ws.Cells.Item(132, 11).Value = 3096  # K132: 3105 -> 3096
$ws.Cells.Item(132, 13).Value = -566  # M132: -575 -> -566

# WVR row 17
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(17, 8).Value = 4000  # H17: 0 -> 4000
$ws.Cells.Item(17, 9).Value = 4000  # I17: 0 -> 4000
$ws.Cells.Item(17, 11).Value = 4000  # K17: 0 -> 4000
$ws.Cells.Item(17, 13).Value = -3828  # M17: None -> -3828

# WVR row 29
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(29, 8).Value = 8582.727999999999  # H29: 8950.666999999999 -> 8582.727999999999
$ws.Cells.Item(29, 9).Value = 441  # I29: 1582.5454 -> 441
$ws.Cells.Item(29, 11).Value = 441  # K29: 1582.5454 -> 441
$ws.Cells.Item(29, 13).Value = -151  # M29: -1292.5454 -> -151

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 0  # H54: 38499 -> 0
$ws.Cells.Item(54, 9).Value = 0  # I54: 6999 -> 0
$ws.Cells.Item(54, 10).Value = 0  # J54: 69999 -> 0
$ws.Cells.Item(54, 11).Value = 0  # K54: 6999 -> 0
$ws.Cells.Item(54, 12).ClearContents()  # L54: 69999 -> (removed)
$ws.Cells.Item(54, 13).ClearContents()  # M54: -6479 -> (removed)
$ws.Cells.Item(54, 14).Value = 0  # N54: -71039 -> 0

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1921  # H126: 1964.75 -> 1921
$ws.Cells.Item(126, 9).Value = 1937.5714  # I126: 1987.5714 -> 1937.5714
$ws.Cells.Item(126, 11).Value = 5812.7142  # K126: 5962.7142 -> 5812.7142
$ws.Cells.Item(126, 13).Value = -3342.7142  # M126: -3492.7142 -> -3342.7142
